# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
#
# Two pairs of match rows had their data entered against the wrong
# fixture / team labels. This swaps the full records (id + stats)
# back onto the correct row, and for the second pair also corrects
# which team name (Home/Away) goes with which row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 91 / 92: fully swap all data columns (B:AD) between the two rows.
$row91 = $ws.Range("B91:AD91")
$row92 = $ws.Range("B92:AD92")

$v91 = $row91.Value2
$v92 = $row92.Value2

$row91.Value2 = $v92
$row92.Value2 = $v91

# --- Rows 186 / 187: swap all data columns (B:AD) EXCEPT the
# HomeTeam/AwayTeam columns (E:F), which stay tied to the row
# (the two fixtures keep their own Home/Away teams; only the
# match stats were mixed up).
$row186_bd = $ws.Range("B186:D186")
$row186_gad = $ws.Range("G186:AD186")
$row187_bd = $ws.Range("B187:D187")
$row187_gad = $ws.Range("G187:AD187")

$v186_bd = $row186_bd.Value2
$v186_gad = $row186_gad.Value2
$v187_bd = $row187_bd.Value2
$v187_gad = $row187_gad.Value2

$row186_bd.Value2 = $v187_bd
$row186_gad.Value2 = $v187_gad
$row187_bd.Value2 = $v186_bd
$row187_gad.Value2 = $v186_gad

# --- Fix the Home/Away team labels for rows 186/187: the team
# names were swapped pairwise (Monterrey U23 <-> Unam Pumas U23,
# Mazatlan FC U23 <-> Tijuana U23).
$ws.Range("E186").Value2 = "Unam Pumas U23"
$ws.Range("F186").Value2 = "Tijuana U23"
$ws.Range("E187").Value2 = "Monterrey U23"
$ws.Range("F187").Value2 = "Mazatlan FC U23"
